$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price column cells keep their original text representation
# (trailing zeros, thousands-dot grouping, leading zeros, etc.) instead of
# being auto-coerced into numbers by Excel when a numeric-looking string is
# assigned to a General-formatted cell.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "30.301.46"
$ws.Range("E2").Value = "  +0.09%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.931.95"
$ws.Range("E3").Value = "  +0.23%  "

$ws.Range("E4").Value = "  +0.16%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.7481"
$ws.Range("E5").Value = "  +4.30%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "242.48"
$ws.Range("E6").Value = "  -2.47%  "

$ws.Range("E7").Value = "  +0.18%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "27.71"
$ws.Range("E8").Value = "  -0.24%  "

$ws.Range("E9").Value = "  -0.48%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.07139"
$ws.Range("E10").Value = "  +0.86%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.7799"
$ws.Range("E11").Value = "  -1.39%  "

$ws.Range("E12").Value = "  +0.77%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.937.71"
$ws.Range("E13").Value = "  +0.60%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.389"
$ws.Range("E14").Value = "  -0.05%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "93.08"
$ws.Range("E15").Value = "  -1.89%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "14.51"
$ws.Range("E16").Value = "  -1.09%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "30.299.00"

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "6.055"
$ws.Range("E18").Value = "  +4.96%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "252.10"
$ws.Range("E19").Value = "  -1.68%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.000007940"
$ws.Range("E20").Value = "  -1.41%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "2.172.07"
$ws.Range("E21").Value = "  -0.32%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.9997"
$ws.Range("E22").Value = "  +0.05%  "

$ws.Range("E23").Value = "  +0.17%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "6.673"
$ws.Range("E24").Value = "  -2.28%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "9.533"
$ws.Range("E25").Value = "  -0.03%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "164.48"
$ws.Range("E26").Value = "  -0.32%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "19.11"
$ws.Range("E27").Value = "  -0.02%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.1294"
$ws.Range("E28").Value = "  +2.43%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.187"
$ws.Range("E29").Value = "  -3.14%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.368"
$ws.Range("E30").Value = "  +0.61%  "

$ws.Range("E31").Value = "  +1.40%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.404"
$ws.Range("E32").Value = "  +0.07%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.138"
$ws.Range("E33").Value = "  +0.43%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.05222"
$ws.Range("E34").Value = "  +1.61%  "

$ws.Range("E35").Value = "  +4.02%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.7554"
$ws.Range("E36").Value = "  +1.40%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.786"
$ws.Range("E37").Value = "  +0.97%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.01952"
$ws.Range("E38").Value = "  -0.31%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.799"
$ws.Range("E39").Value = "  +0.01%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "78.49"
$ws.Range("E40").Value = "  +0.54%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "6.484"
$ws.Range("E41").Value = "  +1.85%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.4514"
$ws.Range("E42").Value = "  +0.31%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.975"
$ws.Range("E43").Value = "  -0.61%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.001"
$ws.Range("E44").Value = "  +0.21%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.8401"
$ws.Range("E45").Value = "  -0.70%  "

$ws.Range("B46").Value = "Aptos"
$ws.Range("C46").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "7.687"
$ws.Range("E46").Value = "  +3.64%  "

$ws.Range("B47").Value = "EnergySwap"
$ws.Range("C47").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "9.984"
$ws.Range("E47").Value = "  +2.45%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "101.75"
$ws.Range("E48").Value = "  +1.27%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "37.84"
$ws.Range("E49").Value = "  +3.42%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.1212"
$ws.Range("E50").Value = "  +6.75%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "953.06"
$ws.Range("E51").Value = "  +3.05%  "
